$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.79505811484151
$ws.Range("C2").Value = 9.25342911341879
$ws.Range("D2").Value = 4.057763500715835
$ws.Range("E2").Value = 11.54444127637297
$ws.Range("F2").Value = 21.2116235354226
$ws.Range("M2").Value = 14.53772324066791
$ws.Range("O2").Value = 18.84733671906793
$ws.Range("B3").Value = 12.14100485638194
$ws.Range("C3").Value = 8.823131878925045
$ws.Range("D3").Value = 4.015209686612613
$ws.Range("E3").Value = 11.47080857032249
$ws.Range("F3").Value = 21.20815380122405
$ws.Range("M3").Value = 14.22207651298686
$ws.Range("O3").Value = 18.91520893882564
$ws.Range("B4").Value = 11.72132454061441
$ws.Range("C4").Value = 8.546677294980002
$ws.Range("D4").Value = 3.988735954430692
$ws.Range("E4").Value = 11.4304968455905
$ws.Range("F4").Value = 21.21499624123001
$ws.Range("M4").Value = 14.02759131164497
$ws.Range("O4").Value = 18.96380518993146
$ws.Range("B5").Value = 11.54593226385322
$ws.Range("C5").Value = 8.431038499508889
$ws.Range("D5").Value = 3.977867639515526
$ws.Range("E5").Value = 11.41531406874008
$ws.Range("F5").Value = 21.22003410938305
$ws.Range("M5").Value = 13.94828026664586
$ws.Range("O5").Value = 18.98533822511556
$ws.Range("B6").Value = 11.51655050331582
$ws.Range("C6").Value = 8.411659923897888
$ws.Range("D6").Value = 3.976058329993121
$ws.Range("E6").Value = 11.41286846491391
$ws.Range("F6").Value = 21.22100626224109
$ws.Range("M6").Value = 13.9351106765938
$ws.Range("O6").Value = 18.98901790344992
$ws.Range("B7").Value = 11.71897656887422
$ws.Range("C7").Value = 8.545129677889705
$ws.Range("D7").Value = 3.988589695610226
$ws.Range("E7").Value = 11.43028703146088
$ws.Range("F7").Value = 21.2150550863212
$ws.Range("M7").Value = 14.02652177273795
$ws.Range("O7").Value = 18.96408860366136
$ws.Range("B8").Value = 12.57339304747119
$ws.Range("C8").Value = 9.107659737878992
$ws.Range("D8").Value = 4.04316672844626
$ws.Range("E8").Value = 11.51804437614496
$ws.Range("F8").Value = 21.20856152742327
$ws.Range("M8").Value = 14.42909576005058
$ws.Range("O8").Value = 18.86929584529748
$ws.Range("B9").Value = 14.09897259975207
$ws.Range("C9").Value = 10.10999521382393
$ws.Range("D9").Value = 4.147156059211881
$ws.Range("E9").Value = 11.72831648210889
$ws.Range("F9").Value = 21.26720711001023
$ws.Range("M9").Value = 15.20837179561114
$ws.Range("O9").Value = 18.73880019948344
$ws.Range("B10").Value = 15.12158952834521
$ws.Range("C10").Value = 10.78109782053151
$ws.Range("D10").Value = 4.22132077808594
$ws.Range("E10").Value = 11.90497093129738
$ws.Range("F10").Value = 21.35388469130676
$ws.Range("M10").Value = 15.76859440031216
$ws.Range("O10").Value = 18.67728926618781
$ws.Range("B11").Value = 15.56444518437468
$ws.Range("C11").Value = 11.07162116396645
$ws.Range("D11").Value = 4.254495616632835
$ws.Range("E11").Value = 11.98986400231085
$ws.Range("F11").Value = 21.40274170417257
$ws.Range("M11").Value = 16.01959224878223
$ws.Range("O11").Value = 18.65688088750017
$ws.Range("B12").Value = 15.72886609835306
$ws.Range("C12").Value = 11.17947343652028
$ws.Range("D12").Value = 4.266970654267781
$ws.Range("E12").Value = 12.02263562002407
$ws.Range("F12").Value = 21.42259045991321
$ws.Range("M12").Value = 16.11399380570401
$ws.Range("O12").Value = 18.65024980406745
$ws.Range("B13").Value = 15.69360186619621
$ws.Range("C13").Value = 11.15634225287202
$ws.Range("D13").Value = 4.264287928499021
$ws.Range("E13").Value = 12.01555034473393
$ws.Range("F13").Value = 21.4182558923073
$ws.Range("M13").Value = 16.09369286934603
$ws.Range("O13").Value = 18.65162901161905
$ws.Range("B14").Value = 15.57803824155644
$ws.Range("C14").Value = 11.08053777710716
$ws.Range("D14").Value = 4.255523737236902
$ws.Range("E14").Value = 11.99254779010716
$ws.Range("F14").Value = 21.40434771901517
$ws.Range("M14").Value = 16.02737217223863
$ws.Range("O14").Value = 18.65631331285573
$ws.Range("B15").Value = 15.50682341952733
$ws.Range("C15").Value = 11.03382266644391
$ws.Range("D15").Value = 4.250143822745506
$ws.Range("E15").Value = 11.97853854586353
$ws.Range("F15").Value = 21.39600378582209
$ws.Range("M15").Value = 15.98666208653555
$ws.Range("O15").Value = 18.6593256867897
$ws.Range("B16").Value = 15.09219235826528
$ws.Range("C16").Value = 10.76181078693402
$ws.Range("D16").Value = 4.219140826797373
$ws.Range("E16").Value = 11.89951187626085
$ws.Range("F16").Value = 21.35088085006237
$ws.Range("M16").Value = 15.75210559929972
$ws.Range("O16").Value = 18.67877610687568
$ws.Range("B17").Value = 14.83205648302689
$ws.Range("C17").Value = 10.5911282008131
$ws.Range("D17").Value = 4.199972533661886
$ws.Range("E17").Value = 11.85217351612605
$ws.Range("F17").Value = 21.32560908828809
$ws.Range("M17").Value = 15.60715907849706
$ws.Range("O17").Value = 18.69265403353904
$ws.Range("B18").Value = 14.68033370125051
$ws.Range("C18").Value = 10.49156821245636
$ws.Range("D18").Value = 4.188894822603188
$ws.Range("E18").Value = 11.82537407436008
$ws.Range("F18").Value = 21.31196139605505
$ws.Range("M18").Value = 15.52343239917072
$ws.Range("O18").Value = 18.70134864789435
$ws.Range("B19").Value = 14.62860463990202
$ws.Range("C19").Value = 10.45762186984085
$ws.Range("D19").Value = 4.185135266915792
$ws.Range("E19").Value = 11.81637461187917
$ws.Range("F19").Value = 21.30749322614092
$ws.Range("M19").Value = 15.4950255334012
$ws.Range("O19").Value = 18.70441458932848
$ws.Range("B20").Value = 14.85996620816981
$ws.Range("C20").Value = 10.60944160363785
$ws.Range("D20").Value = 4.202018525024258
$ws.Range("E20").Value = 11.85716862668395
$ws.Range("F20").Value = 21.32820745896985
$ws.Range("M20").Value = 15.62262654735233
$ws.Range("O20").Value = 18.69110291260539
$ws.Range("B21").Value = 15.61207153456377
$ws.Range("C21").Value = 11.10286235973178
$ws.Range("D21").Value = 4.258100420375604
$ws.Range("E21").Value = 11.99928747600851
$ws.Range("F21").Value = 21.40839638792788
$ws.Range("M21").Value = 16.0468704017787
$ws.Range("O21").Value = 18.65490758399685
$ws.Range("B22").Value = 16.08448177233265
$ws.Range("C22").Value = 11.41272268606469
$ws.Range("D22").Value = 4.294239933918087
$ws.Range("E22").Value = 12.09579638627977
$ws.Range("F22").Value = 21.46865414012546
$ws.Range("M22").Value = 16.32033522414047
$ws.Range("O22").Value = 18.63764988466365
$ws.Range("B23").Value = 15.83411703479725
$ws.Range("C23").Value = 11.24851018139966
$ws.Range("D23").Value = 4.275000683146411
$ws.Range("E23").Value = 12.04396540770485
$ws.Range("F23").Value = 21.43577851316298
$ws.Range("M23").Value = 16.17475873923047
$ws.Range("O23").Value = 18.64627276864222
$ws.Range("B24").Value = 14.84735496661251
$ws.Range("C24").Value = 10.60116657250661
$ws.Range("D24").Value = 4.201093711234425
$ws.Range("E24").Value = 11.85490903931289
$ws.Range("F24").Value = 21.32702998952215
$ws.Range("M24").Value = 15.61563492792555
$ws.Range("O24").Value = 18.69180194461335
$ws.Range("B25").Value = 13.70310366021958
$ws.Range("C25").Value = 9.850071848488298
$ws.Range("D25").Value = 4.119389504970663
$ws.Range("E25").Value = 11.6674496380111
$ws.Range("F25").Value = 21.24368185889164
$ws.Range("M25").Value = 14.9992964751391
$ws.Range("O25").Value = 18.76810836994707
